# Fixing errors in example upload files.
#
# 1) "Service Contacts" sheet: give column A an explicit width and move the
#    selection to D3.
# 2) "Practitioners" sheet: give columns A, C and F explicit widths, add a
#    missing practitioner record (row 6) and move the selection to column G.
# 3) Restore the originally active sheet ("Episodes") so the workbook's
#    active tab is unchanged.

$wb = $excel.ActiveWorkbook

$wsServiceContacts = $wb.Worksheets.Item("Service Contacts")
$wsPractitioners    = $wb.Worksheets.Item("Practitioners")
$wsEpisodes         = $wb.Worksheets.Item("Episodes")

# --- Service Contacts: column width + selection ---------------------------
# ColumnWidth is specified in "characters"; the engine adds ~0.8333 (5/6)
# characters of padding when it stores the width, so back that off here to
# land on a stored width of 14.5.
$wsServiceContacts.Columns.Item(1).ColumnWidth = 13.666666666666666

$wsServiceContacts.Range("D3").Select()

# --- Practitioners: column widths -------------------------------------------
$wsPractitioners.Columns.Item(1).ColumnWidth = 13.833333333333332
$wsPractitioners.Columns.Item(3).ColumnWidth = 12.166666666666666
$wsPractitioners.Columns.Item(6).ColumnWidth = 12.0

# --- Practitioners: new row 6 (missing practitioner record) ---------------
$wsPractitioners.Range("A6").Value = "PHN999:NFP02"
$wsPractitioners.Range("B6").Value = "P01"
$wsPractitioners.Range("C6").Value = 8
$wsPractitioners.Range("D6").Value = 1
$wsPractitioners.Range("E6").Value = 1973
$wsPractitioners.Range("F6").Value = 2
$wsPractitioners.Range("G6").Value = 1
$wsPractitioners.Range("H6").Value = 1
$wsPractitioners.Range("I6").Value = "tag1"

# --- Practitioners: selection ----------------------------------------------
$wsPractitioners.Range("G1:G1048576").Select()

# --- Restore original active sheet -----------------------------------------
$wsEpisodes.Activate()
